$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.016.33'
$ws.Range('E2').Value = '  +0.61%  '

$ws.Range('D3').Value = '2.370.03'
$ws.Range('E3').Value = '  -0.06%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.25'
$ws.Range('E5').Value = '  -3.52%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.11'
$ws.Range('E6').Value = '  +8.92%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.638'
$ws.Range('E7').Value = '  +0.23%  '

$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.618'
$ws.Range('E9').Value = '  -0.18%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.90'
$ws.Range('E10').Value = '  +1.52%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0930'
$ws.Range('E11').Value = '  +0.88%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.61'
$ws.Range('E12').Value = '  +2.16%  '

$ws.Range('E13').Value = '  -0.85%  '

$ws.Range('E14').Value = '  +1.50%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.04'
$ws.Range('E15').Value = '  -1.88%  '

$ws.Range('D16').Value = '2.728.79'
$ws.Range('E16').Value = '  +0.05%  '

$ws.Range('D17').Value = '2.375.08'
$ws.Range('E17').Value = '  +0.31%  '

$ws.Range('D18').Value = '43.045.94'
$ws.Range('E18').Value = '  +0.78%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.69'
$ws.Range('E19').Value = '  -0.15%  '

$ws.Range('E20').Value = '  +0.57%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '76.62'
$ws.Range('E21').Value = '  +1.88%  '

$ws.Range('E22').Value = '  -2.85%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '270.62'
$ws.Range('E23').Value = '  -1.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.34'
$ws.Range('E24').Value = '  +0.81%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.59'
$ws.Range('E25').Value = '  -2.04%  '

$ws.Range('E26').Value = '  +0.43%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.44'
$ws.Range('E27').Value = '  -0.46%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.58'

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.26'
$ws.Range('E29').Value = '  +2.06%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.87'
$ws.Range('E30').Value = '  +3.71%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '167.77'
$ws.Range('E31').Value = '  -3.29%  '

$ws.Range('E32').Value = '  +1.09%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.14'
$ws.Range('E33').Value = '  +4.31%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.93'
$ws.Range('E34').Value = '  -5.82%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.122'
$ws.Range('E35').Value = '  +15.78%  '

$ws.Range('E36').Value = '  -0.03%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.71'
$ws.Range('E37').Value = '  +2.25%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0361'
$ws.Range('E38').Value = '  +0.93%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.84'
$ws.Range('E39').Value = '  -1.46%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.69'
$ws.Range('E40').Value = '  -6.77%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '105.43'
$ws.Range('E41').Value = '  +9.21%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.52'
$ws.Range('E42').Value = '  +0.13%  '

$ws.Range('E43').Value = '  +4.83%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.68'
$ws.Range('E44').Value = '  +4.15%  '

$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.01%  '

$ws.Range('B46').Value = 'Celestia'
$ws.Range('C46').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.66'
$ws.Range('E46').Value = '  +5.96%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '114.69'
$ws.Range('E47').Value = '  -1.65%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '80.46'
$ws.Range('E48').Value = '  +18.22%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.57'
$ws.Range('E49').Value = '  +2.08%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.16'
$ws.Range('E50').Value = '  +2.06%  '

$ws.Range('E51').Value = '  +2.85%  '
